# test-23.xlsx : "mse (3Ysum)" sheet - fill in K5:K54 and P5:P54 forecast values
# (extrapol 2.0 & hybrid3 2.0 (mse)), update the view/selection and column widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mse (3Ysum)")   # sheetId 6 / rId6 / sheet6.xml

# --- Values for column K (rows 5-54) ---
$kVals = @(
    97505.714762790682,97391.868350061181,107150.3802609547,82607.412100734407,92551.962875152996,
    82721.72344369645,102206.4837348837,87344.492795348822,83742.228321664603,104475.83710134641,
    87194.811403304775,89286.856252509169,89107.413994981631,91393.715976744177,103321.5476397797,
    80308.809951529984,72118.470334761325,103609.1791730722,92727.942221664634,87004.955337821273,
    84090.501052264372,84970.028297184806,93954.534851652366,88213.385052753962,97544.262790697685,
    87958.745657282736,96007.633729498164,74649.246682496916,90396.302598408802,84908.145806609522,
    117772.2593356181,78622.491770624227,106037.3946949816,84587.503971113823,91582.478769033012,
    84789.030697552022,86562.352121664619,100049.0954356181,106457.3437903305,91419.526717625457,
    91011.98656646267,106835.6968168911,90217.921144430831,87243.938466095467,94457.17440697673,
    91582.329269767433,87569.278293023235,84722.701322154215,78730.575311015898,101609.2402314565
)

# --- Values for column P (rows 5-54) ---
$pVals = @(
    97249.526480905741,78548.746614810269,81115.498075887386,80880.733705140752,88387.025123255822,
    83042.990164871473,82050.116771970614,84624.566461444309,104300.0243080783,80667.723119339047,
    96237.263360709898,86121.93094039167,83874.836241248457,85578.378965850672,87821.554435618105,
    89165.804376621774,94742.860203304779,72041.439153733168,78759.266331334133,95574.872751897172,
    77014.386386291299,93054.168052631576,94663.463553121182,80138.179235006115,89454.138774785781,
    72150.426653610775,90690.874816768643,85887.11223341494,68118.509080905758,98749.744679559371,
    95859.72477123623,71606.475465728276,89715.730296572816,90647.304252019589,67584.314030477341,
    81921.458703182361,89700.619584577697,88660.380602203179,85589.857345165234,77261.478046389224,
    93607.082694369645,87619.254756303533,101372.9401408813,94505.361277845761,84143.758284088122,
    96352.754174296191,88124.540556058739,81590.952621664619,98547.61101664626,86863.704526560585
)

$n = $kVals.Count

$kArr = New-Object 'object[,]' $n,1
$pArr = New-Object 'object[,]' $n,1
for ($i = 0; $i -lt $n; $i++) {
    $kArr[$i,0] = $kVals[$i]
    $pArr[$i,0] = $pVals[$i]
}

$ws.Range("K5:K54").Value = $kArr
$ws.Range("P5:P54").Value = $pArr

# --- Column widths (new custom widths for columns K=11 and P=16) ---
$ws.Columns.Item(11).ColumnWidth = 11.166666666666666
$ws.Columns.Item(16).ColumnWidth = 9.307291666666666

# --- View / selection changes ---
$ws.Activate()
$ws.Range("V45").Select() | Out-Null

$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 28
